# 1st changes of mifos to finflux
#
# The "Repayment Schedule" sheet gets a new blank column inserted before the
# existing "Late" column (column N), pushing the "Late" and "Outstanding"
# columns one place to the right (N->O, O->P, P->Q). The previously active
# sheet ("Transactions") loses focus and the "Repayment Schedule" sheet
# becomes the active/selected tab instead.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at position N (14th column), shifting the
# "Late" (N) and the blank spacer + "Outstanding" (O, P) columns right.
$ws.Columns.Item(14).Insert()

# Make the "Repayment Schedule" sheet the active tab/sheet and move the
# selection to S4, matching the new sheetView state.
$ws.Select()
$ws.Range("S4").Select()
